# Yearly balance-sheet refresh for سخوز-سیمان خوزستان (Pouya Finance export).
# The 5-year rolling window (columns D:H) shifts one fiscal year to the left
# (D<-E, E<-F, F<-G, G<-H) and the newest fiscal year (1401/12) is appended
# into column H for every data row, including the corresponding period /
# publish-date headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: fiscal period headers (shift left, add FY1401/12) ---
$ws.Cells.Item(8, 4).Value = "12 ماهه منتهی به 1397/12"
$ws.Cells.Item(8, 5).Value = "12 ماهه منتهی به 1398/12"
$ws.Cells.Item(8, 6).Value = "12 ماهه منتهی به 1399/12"
$ws.Cells.Item(8, 7).Value = "12 ماهه منتهی به 1400/12"
$ws.Cells.Item(8, 8).Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish-date headers ---
$ws.Cells.Item(9, 4).Value = "1399-03-21 (10)"
$ws.Cells.Item(9, 5).Value = "1400-02-31 (10)"
$ws.Cells.Item(9, 6).Value = "1401-02-31 (11)"
$ws.Cells.Item(9, 7).Value = "1402-02-25 (10)"
$ws.Cells.Item(9, 8).Value = "1402-02-25 (2)"

# --- Data rows 12-58: shift each 5-year window left by one year and append new FY ---
$ws.Cells.Item(12, 4).Value = 202117
$ws.Cells.Item(12, 5).Value = 1090655
$ws.Cells.Item(12, 6).Value = 812517
$ws.Cells.Item(12, 7).Value = 1664891
$ws.Cells.Item(12, 8).Value = 2657283
$ws.Cells.Item(13, 4).Value = 52446
$ws.Cells.Item(13, 5).Value = 47210
$ws.Cells.Item(13, 6).Value = 1367373
$ws.Cells.Item(13, 7).Value = 2067173
$ws.Cells.Item(13, 8).Value = 5362173
$ws.Cells.Item(14, 4).Value = 632691
$ws.Cells.Item(14, 5).Value = 859938
$ws.Cells.Item(14, 6).Value = 1238474
$ws.Cells.Item(14, 7).Value = 1299988
$ws.Cells.Item(14, 8).Value = 2830471
$ws.Cells.Item(15, 4).Value = 1095555
$ws.Cells.Item(15, 5).Value = 684787
$ws.Cells.Item(15, 6).Value = 1009729
$ws.Cells.Item(15, 7).Value = 2041796
$ws.Cells.Item(15, 8).Value = 3214507
$ws.Cells.Item(16, 4).Value = 403631
$ws.Cells.Item(16, 5).Value = 456675
$ws.Cells.Item(16, 6).Value = 395222
$ws.Cells.Item(16, 7).Value = 751845
$ws.Cells.Item(16, 8).Value = 542690
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 3154
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(18, 4).Value = 2386440
$ws.Cells.Item(18, 5).Value = 3142419
$ws.Cells.Item(18, 6).Value = 4823315
$ws.Cells.Item(18, 7).Value = 7825693
$ws.Cells.Item(18, 8).Value = 14607124
$ws.Cells.Item(19, 4).Value = 5479
$ws.Cells.Item(19, 5).Value = 5996
$ws.Cells.Item(19, 6).Value = 10064
$ws.Cells.Item(19, 7).Value = 11778
$ws.Cells.Item(19, 8).Value = 29738
$ws.Cells.Item(20, 4).Value = 625915
$ws.Cells.Item(20, 5).Value = 766722
$ws.Cells.Item(20, 6).Value = 725758
$ws.Cells.Item(20, 7).Value = 819464
$ws.Cells.Item(20, 8).Value = 931763
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(22, 4).Value = 1314471
$ws.Cells.Item(22, 5).Value = 1202619
$ws.Cells.Item(22, 6).Value = 934792
$ws.Cells.Item(22, 7).Value = 1350535
$ws.Cells.Item(22, 8).Value = 1435286
$ws.Cells.Item(23, 4).Value = 66771
$ws.Cells.Item(23, 5).Value = 66760
$ws.Cells.Item(23, 6).Value = 93030
$ws.Cells.Item(23, 7).Value = 93030
$ws.Cells.Item(23, 8).Value = 93030
$ws.Cells.Item(24, 4).Value = "-"
$ws.Cells.Item(24, 5).Value = "-"
$ws.Cells.Item(24, 6).Value = "-"
$ws.Cells.Item(24, 7).Value = "-"
$ws.Cells.Item(24, 8).Value = "-"
$ws.Cells.Item(25, 4).Value = 30722
$ws.Cells.Item(25, 5).Value = 49667
$ws.Cells.Item(25, 6).Value = 33376
$ws.Cells.Item(25, 7).Value = 37781
$ws.Cells.Item(25, 8).Value = 30460
$ws.Cells.Item(26, 4).Value = 2043358
$ws.Cells.Item(26, 5).Value = 2091764
$ws.Cells.Item(26, 6).Value = 1797020
$ws.Cells.Item(26, 7).Value = 2312588
$ws.Cells.Item(26, 8).Value = 2520277
$ws.Cells.Item(27, 4).Value = 4429798
$ws.Cells.Item(27, 5).Value = 5234183
$ws.Cells.Item(27, 6).Value = 6620335
$ws.Cells.Item(27, 7).Value = 10138281
$ws.Cells.Item(27, 8).Value = 17127401
$ws.Cells.Item(29, 4).Value = 698702
$ws.Cells.Item(29, 5).Value = 1450924
$ws.Cells.Item(29, 6).Value = 820268
$ws.Cells.Item(29, 7).Value = 1571666
$ws.Cells.Item(29, 8).Value = 2261831
$ws.Cells.Item(30, 4).Value = "-"
$ws.Cells.Item(30, 5).Value = "-"
$ws.Cells.Item(30, 6).Value = "-"
$ws.Cells.Item(30, 7).Value = "-"
$ws.Cells.Item(30, 8).Value = "-"
$ws.Cells.Item(31, 4).Value = 95258
$ws.Cells.Item(31, 5).Value = 162316
$ws.Cells.Item(31, 6).Value = 591696
$ws.Cells.Item(31, 7).Value = 1007990
$ws.Cells.Item(31, 8).Value = 1017482
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = 58271
$ws.Cells.Item(32, 6).Value = 194106
$ws.Cells.Item(32, 7).Value = 257899
$ws.Cells.Item(32, 8).Value = 596108
$ws.Cells.Item(33, 4).Value = 749776
$ws.Cells.Item(33, 5).Value = 153045
$ws.Cells.Item(33, 6).Value = 107917
$ws.Cells.Item(33, 7).Value = 226929
$ws.Cells.Item(33, 8).Value = 250332
$ws.Cells.Item(34, 4).Value = 1048293
$ws.Cells.Item(34, 5).Value = 553718
$ws.Cells.Item(34, 6).Value = 86865
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 457667
$ws.Cells.Item(35, 4).Value = 0
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(36, 4).Value = 0
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(37, 4).Value = 2592029
$ws.Cells.Item(37, 5).Value = 2378274
$ws.Cells.Item(37, 6).Value = 1800852
$ws.Cells.Item(37, 7).Value = 3064484
$ws.Cells.Item(37, 8).Value = 4583420
$ws.Cells.Item(38, 4).Value = 0
$ws.Cells.Item(38, 5).Value = 0
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(39, 4).Value = "-"
$ws.Cells.Item(39, 5).Value = "-"
$ws.Cells.Item(39, 6).Value = "-"
$ws.Cells.Item(39, 7).Value = "-"
$ws.Cells.Item(39, 8).Value = "-"
$ws.Cells.Item(40, 4).Value = 0
$ws.Cells.Item(40, 5).Value = 0
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(41, 4).Value = 85888
$ws.Cells.Item(41, 5).Value = 89430
$ws.Cells.Item(41, 6).Value = 121249
$ws.Cells.Item(41, 7).Value = 181100
$ws.Cells.Item(41, 8).Value = 309809
$ws.Cells.Item(42, 4).Value = 85888
$ws.Cells.Item(42, 5).Value = 89430
$ws.Cells.Item(42, 6).Value = 121249
$ws.Cells.Item(42, 7).Value = 181100
$ws.Cells.Item(42, 8).Value = 309809
$ws.Cells.Item(43, 4).Value = 2677917
$ws.Cells.Item(43, 5).Value = 2467704
$ws.Cells.Item(43, 6).Value = 1922101
$ws.Cells.Item(43, 7).Value = 3245584
$ws.Cells.Item(43, 8).Value = 4893229
$ws.Cells.Item(45, 4).Value = 650000
$ws.Cells.Item(45, 5).Value = 650000
$ws.Cells.Item(45, 6).Value = 1400000
$ws.Cells.Item(45, 7).Value = 1400000
$ws.Cells.Item(45, 8).Value = 1400000
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 0
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 646136
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(48, 4).Value = 0
$ws.Cells.Item(48, 5).Value = 0
$ws.Cells.Item(48, 6).Value = -84358
$ws.Cells.Item(48, 7).Value = -59888
$ws.Cells.Item(48, 8).Value = -8028
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = 0
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 19131
$ws.Cells.Item(49, 8).Value = 81477
$ws.Cells.Item(50, 4).Value = 65000
$ws.Cells.Item(50, 5).Value = 65000
$ws.Cells.Item(50, 6).Value = 140000
$ws.Cells.Item(50, 7).Value = 140000
$ws.Cells.Item(50, 8).Value = 140000
$ws.Cells.Item(51, 4).Value = 529
$ws.Cells.Item(51, 5).Value = 529
$ws.Cells.Item(51, 6).Value = 529
$ws.Cells.Item(51, 7).Value = 529
$ws.Cells.Item(51, 8).Value = 529
$ws.Cells.Item(52, 4).Value = "-"
$ws.Cells.Item(52, 5).Value = "-"
$ws.Cells.Item(52, 6).Value = "-"
$ws.Cells.Item(52, 7).Value = "-"
$ws.Cells.Item(52, 8).Value = "-"
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(54, 4).Value = "-"
$ws.Cells.Item(54, 5).Value = "-"
$ws.Cells.Item(54, 6).Value = "-"
$ws.Cells.Item(54, 7).Value = "-"
$ws.Cells.Item(54, 8).Value = "-"
$ws.Cells.Item(55, 4).Value = 0
$ws.Cells.Item(55, 5).Value = 0
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(56, 4).Value = 1036352
$ws.Cells.Item(56, 5).Value = 1404814
$ws.Cells.Item(56, 6).Value = 3242063
$ws.Cells.Item(56, 7).Value = 5392925
$ws.Cells.Item(56, 8).Value = 10620194
$ws.Cells.Item(57, 4).Value = 1751881
$ws.Cells.Item(57, 5).Value = 2766479
$ws.Cells.Item(57, 6).Value = 4698234
$ws.Cells.Item(57, 7).Value = 6892697
$ws.Cells.Item(57, 8).Value = 12234172
$ws.Cells.Item(58, 4).Value = 4429798
$ws.Cells.Item(58, 5).Value = 5234183
$ws.Cells.Item(58, 6).Value = 6620335
$ws.Cells.Item(58, 7).Value = 10138281
$ws.Cells.Item(58, 8).Value = 17127401
